$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 2286398.04
$ws.Range("C7").Value = -48.5402241813385
$ws.Range("D7").Value = 2314
$ws.Range("E7").Value = 2314
$ws.Range("F7").Value = 988.0717545375973
$ws.Range("G7").Value = 5.321304354875056
